$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the record data between row 3 and row 4 for columns A, B, E, F, G, H, Q, R.
# (C, D, I, K, P, S and the remaining columns stay identical between the two rows,
# so they do not need to be touched.)
$cols = @("A", "B", "E", "F", "G", "H", "Q", "R")

foreach ($col in $cols) {
    $rangeRow3 = $ws.Range("$col`3")
    $rangeRow4 = $ws.Range("$col`4")

    $valueRow3 = $rangeRow3.Value2
    $valueRow4 = $rangeRow4.Value2

    $rangeRow3.Value2 = $valueRow4
    $rangeRow4.Value2 = $valueRow3
}
